# BIO-iTC-CommentsMatrix.xlsx - publish toolbox PR1 for review
# - Extend the "DOCUMENT" dropdown list on the Instructions sheet with
#   Eye / Face / Finger / Vein toolbox entries (H12:H15) and grow the
#   backing Table1 (and its AutoFilter) to cover H7:H15.
# - Update the instructional text in B3 to document the four new codes
#   and grow the row to fit the extra lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# New dropdown entries, displayed in alphabetical order beneath the
# existing cPP/CFG/SD/TB rows.
$ws.Range("H12").Value = "Eye"
$ws.Range("H13").Value = "Face"
$ws.Range("H14").Value = "Finger"
$ws.Range("H15").Value = "Vein"

# Grow the DOCUMENT table (Table1) so its range / AutoFilter cover the
# newly added rows - this keeps the G1 dropdown list in sync.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("H7:H15"))

# Update the explanatory note for the document-type field, and resize
# the row so the added lines are fully visible.
$newNote = "The specific document the comments are for should be selected in G1. The fields are:" + "`n" + `
  "- cPP - for the PP-Module" + "`n" + `
  "- CFG - for the PP-Configuration" + "`n" + `
  "- SD - for the Supporting Document" + "`n" + `
  "- TB - for the PAD Toolbox overview" + "`n" + `
  "- Eye - for the Eye Toolbox" + "`n" + `
  "- Face - for the Face Toolbox" + "`n" + `
  "- Finger - for the Fingerprint Toolbox" + "`n" + `
  "- Vein - for the Vein Toolbox"

$ws.Range("B3").Value = $newNote
$ws.Rows.Item(3).RowHeight = 141.75

Write-Output "Instructions sheet updated for Eye/Face/Finger/Vein toolboxes"
